$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'299.22"
$ws.Range('E2').Value = "'1.60%"
$ws.Range('E3').Value = "'-0.09%"
$ws.Range('D4').Value = "'5.136"
$ws.Range('E4').Value = "'0.50%"
$ws.Range('D5').Value = "'0.08073"
$ws.Range('E5').Value = "'9.70%"
$ws.Range('D6').Value = "'2.671"
$ws.Range('E6').Value = "'62.88%"
$ws.Range('D7').Value = "'7.852"
$ws.Range('E7').Value = "'2.31%"
$ws.Range('E8').Value = "'2.05%"
$ws.Range('D9').Value = "'0.9094"
$ws.Range('E9').Value = "'-1.12%"
$ws.Range('D10').Value = "'0.1729"
$ws.Range('E10').Value = "'3.41%"
$ws.Range('D11').Value = "'0.07260"
$ws.Range('E11').Value = "'2.56%"
$ws.Range('D12').Value = "'0.08066"
$ws.Range('E12').Value = "'1.81%"
$ws.Range('D13').Value = "'0.03026"
$ws.Range('E13').Value = "'0.95%"
$ws.Range('D14').Value = "'0.09979"
$ws.Range('E14').Value = "'0.99%"
$ws.Range('D15').Value = "'0.001498"
$ws.Range('E15').Value = "'0.47%"
$ws.Range('D16').Value = "'0.006039"
$ws.Range('E16').Value = "'-2.06%"
$ws.Range('D17').Value = "'3.504"
$ws.Range('E17').Value = "'1.57%"
$ws.Range('E18').Value = "'1.13%"
$ws.Range('D19').Value = "'0.3289"
$ws.Range('E19').Value = "'0.31%"
$ws.Range('E20').Value = "'0.33%"
$ws.Range('D21').Value = "'4.610"
$ws.Range('E21').Value = "'1.20%"
$ws.Range('D22').Value = "'0.1602"
$ws.Range('E22').Value = "'3.24%"
$ws.Range('D23').Value = "'0.04595"
$ws.Range('E23').Value = "'-0.63%"
$ws.Range('D24').Value = "'0.001267"
$ws.Range('E24').Value = "'4.01%"
$ws.Range('D25').Value = "'0.004444"
$ws.Range('E26').Value = "'-9.11%"
$ws.Range('D27').Value = "'0.0003436"
$ws.Range('E27').Value = "'83.14%"
$ws.Range('D39').Value = "'0.01813"
$ws.Range('E39').Value = "'7.40%"
$ws.Range('D40').Value = "'0.04539"
$ws.Range('E40').Value = "'3.18%"
$ws.Range('D41').Value = "'0.007092"
$ws.Range('E41').Value = "'0.18%"
$ws.Range('D42').Value = "'0.1341"
$ws.Range('E42').Value = "'1.27%"
$ws.Range('E43').Value = "'6.30%"
$ws.Range('D44').Value = "'0.01049"
$ws.Range('E44').Value = "'3.20%"
$ws.Range('D45').Value = "'0.00006334"
$ws.Range('E45').Value = "'5.96%"
$ws.Range('E46').Value = "'0.14%"
$ws.Range('B47').Value = 'BOLO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D47').Value = "'0.8206"
$ws.Range('E47').Value = "'-57.22%"
$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').Value = "'0.006410"
$ws.Range('E48').Value = "'-41.73%"
$ws.Range('D49').Value = "'0.00002103"
$ws.Range('E49').Value = "'0.14%"
$ws.Range('D50').Value = "'0.0002003"
$ws.Range('E50').Value = "'0.21%"
